# Apply updated crypto price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.814.47'
$ws.Range('E2').Value = '  +3.16%  '
$ws.Range('D3').Value = '1.865.95'
$ws.Range('E3').Value = '  +2.80%  '
$ws.Range('D4').Value = "'1.041"
$ws.Range('E4').Value = '  +3.45%  '
$ws.Range('D5').Value = "'324.96"
$ws.Range('E5').Value = '  +4.07%  '
$ws.Range('E6').Value = '  +3.12%  '
$ws.Range('D7').Value = "'0.4425"
$ws.Range('E7').Value = '  +2.91%  '
$ws.Range('D8').Value = "'0.3798"
$ws.Range('E8').Value = '  +3.41%  '
$ws.Range('D9').Value = "'0.07477"
$ws.Range('E9').Value = '  +2.96%  '
$ws.Range('D10').Value = "'0.8858"
$ws.Range('E10').Value = '  +2.35%  '
$ws.Range('D11').Value = "'21.75"
$ws.Range('E11').Value = '  +1.96%  '
$ws.Range('D12').Value = '1.888.53'
$ws.Range('E12').Value = '  -12.15%  '
$ws.Range('D13').Value = "'5.562"
$ws.Range('E13').Value = '  +2.85%  '
$ws.Range('D14').Value = "'6.766"
$ws.Range('E14').Value = '  +2.41%  '
$ws.Range('D15').Value = "'0.07233"
$ws.Range('E15').Value = '  +4.21%  '
$ws.Range('D16').Value = "'83.84"
$ws.Range('E16').Value = '  +3.51%  '
$ws.Range('D17').Value = "'1.043"
$ws.Range('E17').Value = '  +3.65%  '
$ws.Range('D18').Value = "'0.000009174"
$ws.Range('E18').Value = '  +3.37%  '
$ws.Range('D19').Value = "'1.036"
$ws.Range('E19').Value = '  +3.08%  '
$ws.Range('D20').Value = "'15.56"
$ws.Range('E20').Value = '  +1.98%  '
$ws.Range('D21').Value = '27.859.79'
$ws.Range('E21').Value = '  +3.19%  '
$ws.Range('D22').Value = "'5.324"
$ws.Range('E22').Value = '  +2.63%  '
$ws.Range('E23').Value = '  +3.29%  '
$ws.Range('D24').Value = "'1.994"
$ws.Range('D25').Value = "'158.64"
$ws.Range('E25').Value = '  +3.06%  '
$ws.Range('E26').Value = '  +2.90%  '
$ws.Range('D27').Value = "'5.344"
$ws.Range('E27').Value = '  +2.33%  '
$ws.Range('D28').Value = "'1.986"
$ws.Range('E28').Value = '  +4.55%  '
$ws.Range('D29').Value = "'117.74"
$ws.Range('E29').Value = '  +2.69%  '
$ws.Range('D30').Value = "'0.09080"
$ws.Range('E30').Value = '  +1.49%  '
$ws.Range('B31').Value = 'HuobiToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D31').Value = "'3.116"
$ws.Range('E31').Value = '  +10.91%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = "'0.7785"
$ws.Range('E32').Value = '  +4.32%  '
$ws.Range('D33').Value = "'1.216"
$ws.Range('E33').Value = '  +2.29%  '
$ws.Range('D34').Value = "'4.577"
$ws.Range('E34').Value = '  +3.60%  '
$ws.Range('E35').Value = '  +3.26%  '
$ws.Range('E36').Value = '  +2.55%  '
$ws.Range('D37').Value = "'0.01996"
$ws.Range('E37').Value = '  +3.71%  '
$ws.Range('D38').Value = "'0.05355"
$ws.Range('E38').Value = '  +2.71%  '
$ws.Range('D39').Value = "'2.881"
$ws.Range('E39').Value = '  +5.32%  '
$ws.Range('D40').Value = "'0.5205"
$ws.Range('E40').Value = '  +2.04%  '
$ws.Range('D41').Value = "'0.1696"
$ws.Range('E41').Value = '  +2.41%  '
$ws.Range('D42').Value = "'6.935"
$ws.Range('E42').Value = '  +7.18%  '
$ws.Range('D43').Value = "'8.687"
$ws.Range('E43').Value = '  +4.53%  '
$ws.Range('D44').Value = "'10.77"
$ws.Range('E44').Value = '  +3.55%  '
$ws.Range('D45').Value = "'109.70"
$ws.Range('E45').Value = '  +2.74%  '
$ws.Range('D46').Value = "'1.728"
$ws.Range('E46').Value = '  +5.09%  '
$ws.Range('D47').Value = "'0.4713"
$ws.Range('E47').Value = '  +2.72%  '
$ws.Range('D48').Value = "'0.06469"
$ws.Range('E48').Value = '  +4.11%  '
$ws.Range('D49').Value = "'1.913"
$ws.Range('E49').Value = '  +3.90%  '
$ws.Range('D50').Value = "'39.89"
$ws.Range('E50').Value = '  +3.86%  '
$ws.Range('E51').Value = '  +2.53%  '
